$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20; all rows from 20 downward shift down by one.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44581
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112021
$ws.Range("G20").Value = "Ají"
$ws.Range("H20").Value = "Americana (o)"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 16000
$ws.Range("L20").Value = 17000
$ws.Range("M20").Value = 16500
$ws.Range("N20").Value = "$/caja 15 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 1100
$ws.Range("Q20").Value = 15
$ws.Range("R20").Value = "Hortaliza"
